$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain plain text, matching the
# original inlineStr string cells (many prices look numeric, e.g. "587.69").
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '65.120.22'
$ws.Range('E2').Value = '  -2.04%  '

$ws.Range('D3').Value = '3.476.40'
$ws.Range('E3').Value = '  -0.91%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '587.69'
$ws.Range('E5').Value = '  -2.78%  '

$ws.Range('D6').Value = '136.93'
$ws.Range('E6').Value = '  -4.51%  '

$ws.Range('D7').Value = '3.476.27'
$ws.Range('E7').Value = '  -0.91%  '

$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('E9').Value = '  -2.81%  '

$ws.Range('E10').Value = '  -5.77%  '

$ws.Range('E11').Value = '  -7.09%  '

$ws.Range('D12').Value = '0.383'
$ws.Range('E12').Value = '  -4.72%  '

$ws.Range('D13').Value = '4.069.17'
$ws.Range('E13').Value = '  -0.67%  '

$ws.Range('D14').Value = '0.0000181'
$ws.Range('E14').Value = '  -6.48%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.484.93'
$ws.Range('E15').Value = '  -0.98%  '

$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').Value = '26.57'
$ws.Range('E16').Value = '  -7.19%  '

$ws.Range('E17').Value = '  -1.38%  '

$ws.Range('D18').Value = '65.105.64'
$ws.Range('E18').Value = '  -1.85%  '

$ws.Range('D19').Value = '9.70'
$ws.Range('E19').Value = '  -8.78%  '

$ws.Range('D20').Value = '5.77'
$ws.Range('E20').Value = '  -5.21%  '

$ws.Range('E21').Value = '  -4.41%  '

$ws.Range('D22').Value = '388.39'
$ws.Range('E22').Value = '  -7.77%  '

$ws.Range('D23').Value = '0.555'
$ws.Range('E23').Value = '  -5.18%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.02%  '

$ws.Range('B25').Value = 'LEO'
$ws.Range('C25').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D25').Value = '5.76'
$ws.Range('E25').Value = '  +0.95%  '

$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '72.45'
$ws.Range('E26').Value = '  -5.49%  '

$ws.Range('D27').Value = '3.618.37'
$ws.Range('E27').Value = '  -0.95%  '

$ws.Range('E28').Value = '  -1.95%  '

$ws.Range('E29').Value = '  -0.09%  '

$ws.Range('D30').Value = '7.36'
$ws.Range('E30').Value = '  -5.00%  '

$ws.Range('D31').Value = '8.12'
$ws.Range('E31').Value = '  -8.85%  '

$ws.Range('E32').Value = '  -9.55%  '

$ws.Range('D33').Value = '3.497.29'
$ws.Range('E33').Value = '  -0.45%  '

$ws.Range('E35').Value = '  -6.79%  '

$ws.Range('D36').Value = '23.03'
$ws.Range('E36').Value = '  -4.48%  '

$ws.Range('D37').Value = '170.77'
$ws.Range('E37').Value = '  -1.52%  '

$ws.Range('E38').Value = '  -9.29%  '

$ws.Range('E39').Value = '  -9.06%  '

$ws.Range('E40').Value = '  -9.16%  '

$ws.Range('D41').Value = '4.72'
$ws.Range('E41').Value = '  -8.62%  '

$ws.Range('D42').Value = '0.0777'
$ws.Range('E42').Value = '  -3.28%  '

$ws.Range('D43').Value = '0.810'
$ws.Range('E43').Value = '  -4.60%  '

$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.09%  '

$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '42.49'
$ws.Range('E45').Value = '  -6.61%  '

$ws.Range('D46').Value = '25.03'
$ws.Range('E46').Value = '  +9.20%  '

$ws.Range('D47').Value = '4.34'
$ws.Range('E47').Value = '  -11.96%  '

$ws.Range('E48').Value = '  +3.91%  '

$ws.Range('E49').Value = '  -8.28%  '

$ws.Range('D50').Value = '6.69'
$ws.Range('E50').Value = '  -4.73%  '

$ws.Range('D51').Value = '2.215.59'
$ws.Range('E51').Value = '  -3.56%  '

# Remove the temporary text format so the cells keep the default style
# (no explicit style index), exactly like the rest of the sheet.
$priceRange.ClearFormats()
